# --- Set zoom on the existing "ReFuelEU" sheet while it is still active ---
$wb = $excel.ActiveWorkbook
$excel.ActiveWindow.Zoom = 135

# --- Add the new "EU Production" sheet, placed after "ReFuelEU" ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "EU Production"
$newSheet.Move($null, $wb.Worksheets.Item("ReFuelEU"))

# re-fetch a fresh handle to the freshly-moved sheet (the old handle can go stale
# after the structural Move edit)
$ws = $wb.Worksheets.Item("EU Production")

# --- Populate cells in the same order the shared-string table records them ---
$ws.Range("A2").Value = "domestic aviation consumption, EU 27 [t(oil equivalent)]"
$ws.Range("A1").Value = "metric"
$ws.Range("B2").Value = "2019 (=pre-COVID)"
$ws.Range("C1").Value = "value"
$ws.Range("A3").Value = "all biofuels production, EU 27 [t(oil equivalent)]"
$ws.Range("E2").Value = "3.1.4 EU-27"
$ws.Range("E3").Value = "3.1.7 Biofuels Production"
$ws.Range("A4").Value = "bio jet fuel, EU 27 [t(oil equivalent)]"
$ws.Range("E1").Value = "source table"
$ws.Range("B1").Value = "year"
$ws.Range("D1").Value = "source"
$ws.Range("E4").Value = "3.1.7 Biofuels Production"

$ws.Range("B3").Value = 2021
$ws.Range("B4").Value = 2021

$ws.Range("C2").Formula = "=6.6*10000000"
$ws.Range("C3").Formula = "=15960*10000"
$ws.Range("C4").Formula = "=91.6*10000"

# --- Hyperlinks (source table column) ---
$url = "https://transport.ec.europa.eu/facts-funding/studies-data/eu-transport-figures-statistical-pocketbook/statistical-pocketbook-2023_en"
$ws.Hyperlinks.Add($ws.Range("D2"), $url)
$ws.Range("D2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D3"), $url)
$ws.Range("D3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D4"), $url)
$ws.Range("D4").Style = "Hyperlink"

# --- Column A width (bestFit 48.5) ---
$ws.Columns.Item(1).ColumnWidth = 47.666666666666664

# --- Activate the new sheet, set its zoom and selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 156
$ws.Range("A12").Select()
